$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header "Save" in H1, matching the style used by the other headers (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Add the "Save" column values (all zero) for rows 2-5
$ws.Range("H2:H5").Value = 0
